$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 27 (2026-02) statistics per commit "Update stats for 2026-02"
$ws.Range("B27").Value = 6540
$ws.Range("C27").Value = 1017
$ws.Range("D27").Value = 6096853
$ws.Range("E27").Value = 932.2405198776759
$ws.Range("F27").Value = 9.915966386554631
$ws.Range("G27").Value = 7.391763463569156
$ws.Range("H27").Value = 25.05064601486573
